$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 173
$ws.Range("A173").Value = "Xasanova Ma'rifat Aminovna"
$ws.Range("B173").Value = "Maktabgacha talim tashkiloti direktori"
$ws.Range("C173").Value = "AB3171624"
$ws.Range("D173").NumberFormat = "@"
$ws.Range("D173").Value = "348"
$ws.Range("E173").Value = "Buxoro viloyati"
$ws.Range("F173").Value = "Kogon tumani"
$ws.Range("G173").NumberFormat = "@"
$ws.Range("G173").Value = "998978608883"
$ws.Range("H173").NumberFormat = "@"
$ws.Range("H173").Value = "10-12-2024"

# New row 174
$ws.Range("A174").Value = "Xashimova Dilnaz Sabitxanovna"
$ws.Range("B174").Value = "Maktabgacha talim tashkiloti tarbiyachisi"
$ws.Range("C174").Value = "AB0959509"
$ws.Range("D174").NumberFormat = "@"
$ws.Range("D174").Value = "349"
$ws.Range("E174").Value = "Toshkent shahri"
$ws.Range("F174").Value = "Yunusobod tumani"
$ws.Range("G174").NumberFormat = "@"
$ws.Range("G174").Value = "998909779885"
$ws.Range("H174").NumberFormat = "@"
$ws.Range("H174").Value = "10-12-2024"
